# Generate Report for Handoff
#
# For every row whose Status is "Ready for handoff" (rows 4-7 on both the
# "zh-cn" and "de-de" language sheets), a new handoff has just been
# generated:
#   - Priority goes from "low" to "ht" (now queued for handoff)
#   - Latest Handoff Datetime is refreshed to the new generation time
#
# zh-cn handoff time: 2016-08-18 12:31:18
# de-de handoff time: 2016-08-18 12:31:23

$wb = $excel.ActiveWorkbook

$newPriority = "ht"

$zhRows = 4..7
$zhHandoffTime = "2016-08-18 12:31:18"

$deRows = 4..7
$deHandoffTime = "2016-08-18 12:31:23"

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $wsZh.Range("E$r").Value = $newPriority
    $wsZh.Range("H$r").Value = $zhHandoffTime
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $wsDe.Range("E$r").Value = $newPriority
    $wsDe.Range("H$r").Value = $deHandoffTime
}

# The "Overview" sheet mirrors the latest handoff-generation timestamp in
# column G ("Latest HO Xliff Generate Date"); it tracks the de-de value.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = $deHandoffTime
}
